$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 353
$ws.Range("I6").Value = 236
$ws.Range("J6").Value = 860
$ws.Range("K6").Value = 708
$ws.Range("L6").Value = 2580
$ws.Range("M6").Value = -596
$ws.Range("N6").Value = -2804

$ws.Range("H40").Value = 1237.4572
$ws.Range("I40").Value = 1200.4783
$ws.Range("J40").Value = 1308.3334
$ws.Range("K40").Value = 1200.4783
$ws.Range("L40").Value = 1308.3334
$ws.Range("M40").Value = -1025.4783
$ws.Range("N40").Value = -1658.3334

$ws.Range("H64").Value = 3019.6155
$ws.Range("J64").Value = 3076.6667
$ws.Range("L64").Value = 3076.6667
$ws.Range("N64").Value = -3572.6667

$ws.Range("H67").Value = 3019.6155
$ws.Range("J67").Value = 3076.6667
$ws.Range("L67").Value = 3076.6667
$ws.Range("N67").Value = -4792.6667

$ws.Range("H76").Value = 3260.9092
$ws.Range("I76").Value = 3210.5557
$ws.Range("J76").Value = 3487.5
$ws.Range("K76").Value = 3210.5557
$ws.Range("L76").Value = 3487.5
$ws.Range("M76").Value = -2895.5557
$ws.Range("N76").Value = -4117.5

$ws.Range("H79").Value = 3260.9092
$ws.Range("I79").Value = 3210.5557
$ws.Range("J79").Value = 3487.5
$ws.Range("K79").Value = 3210.5557
$ws.Range("L79").Value = 3487.5
$ws.Range("M79").Value = -2118.5557
$ws.Range("N79").Value = -5671.5

$ws.Range("H98").Value = 2604.8147
$ws.Range("I98").Value = 1685.4736
$ws.Range("J98").Value = 4788.25
$ws.Range("K98").Value = 1685.4736
$ws.Range("L98").Value = 4788.25
$ws.Range("M98").Value = -187.4736
$ws.Range("N98").Value = -7784.25

$ws.Range("H112").Value = 1565.8667
$ws.Range("I112").Value = 850
$ws.Range("J112").Value = 1617
$ws.Range("K112").Value = 2550
$ws.Range("L112").Value = 4851
$ws.Range("M112").Value = -1442
$ws.Range("N112").Value = -7067

$ws.Range("H122").Value = 2604.8147
$ws.Range("I122").Value = 1685.4736
$ws.Range("J122").Value = 4788.25
$ws.Range("K122").Value = 5056.4208
$ws.Range("L122").Value = 14364.75
$ws.Range("M122").Value = -2606.4208
$ws.Range("N122").Value = -19264.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2719.8857
$ws.Range("I63").Value = 2611.2
$ws.Range("J63").Value = 2864.8
$ws.Range("K63").Value = 2611.2
$ws.Range("L63").Value = 2864.8
$ws.Range("M63").Value = -1925.2
$ws.Range("N63").Value = -4236.8

$ws.Range("H66").Value = 2719.8857
$ws.Range("I66").Value = 2611.2
$ws.Range("J66").Value = 2864.8
$ws.Range("K66").Value = 13056
$ws.Range("L66").Value = 14324
$ws.Range("M66").Value = -9624
$ws.Range("N66").Value = -21188

$ws.Range("H130").Value = 12809.333
$ws.Range("J130").Value = 12809.333
$ws.Range("L130").Value = 12809.333
$ws.Range("N130").Value = -22849.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3678.6365
$ws.Range("I86").Value = 3405.0227
$ws.Range("J86").Value = 4773.091
$ws.Range("K86").Value = 3405.0227
$ws.Range("L86").Value = 4773.091
$ws.Range("M86").Value = -2282.0227
$ws.Range("N86").Value = -7019.091

$ws.Range("H89").Value = 3678.6365
$ws.Range("I89").Value = 3405.0227
$ws.Range("J89").Value = 4773.091
$ws.Range("K89").Value = 17025.1135
$ws.Range("L89").Value = 23865.455
$ws.Range("M89").Value = -11409.1135
$ws.Range("N89").Value = -35097.455

$ws.Range("H134").Value = 1121.3489
$ws.Range("I134").Value = 990.9722
$ws.Range("J134").Value = 1791.8572
$ws.Range("K134").Value = 2972.9166
$ws.Range("L134").Value = 5375.571599999999
$ws.Range("M134").Value = -437.9166
$ws.Range("N134").Value = -10445.5716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2400
$ws.Range("I62").Value = 2400
$ws.Range("J62").Value = 2400
$ws.Range("K62").Value = 2400
$ws.Range("L62").Value = 2400
$ws.Range("M62").Value = -1776
$ws.Range("N62").Value = -3648

$ws.Range("H65").Value = 2400
$ws.Range("I65").Value = 2400
$ws.Range("J65").Value = 2400
$ws.Range("K65").Value = 12000
$ws.Range("L65").Value = 12000
$ws.Range("M65").Value = -8880
$ws.Range("N65").Value = -18240

$ws.Range("H134").Value = 1479.1842
$ws.Range("I134").Value = 1576.9395
$ws.Range("J134").Value = 834
$ws.Range("K134").Value = 4730.818499999999
$ws.Range("L134").Value = 2502
$ws.Range("M134").Value = -2195.818499999999
$ws.Range("N134").Value = -7572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 902.29114
$ws.Range("I131").Value = 519.7
$ws.Range("J131").Value = 957.73914
$ws.Range("K131").Value = 1559.1
$ws.Range("L131").Value = 2873.21742
$ws.Range("M131").Value = 3480.9
$ws.Range("N131").Value = -12953.21742

$ws.Range("H137").Value = 2985.8076
$ws.Range("I137").Value = 1231.4286
$ws.Range("J137").Value = 5032.5835
$ws.Range("K137").Value = 3694.2858
$ws.Range("L137").Value = 15097.7505
$ws.Range("M137").Value = 1405.7142
$ws.Range("N137").Value = -25297.7505

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 33000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 33000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 33000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -33224

$ws.Range("H8").Value = 33000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 33000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 33000
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -33278

$ws.Range("H10").Value = 933.3333
$ws.Range("I10").Value = 933.3333
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 933.3333
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -764.3333
$ws.Range("N10").ClearContents()

$ws.Range("H70").Value = 4138.5
$ws.Range("I70").Value = 4151.067
$ws.Range("J70").Value = 4117.5557
$ws.Range("K70").Value = 4151.067
$ws.Range("L70").Value = 4117.5557
$ws.Range("M70").Value = -3881.067
$ws.Range("N70").Value = -4657.5557

$ws.Range("H73").Value = 4138.5
$ws.Range("I73").Value = 4151.067
$ws.Range("J73").Value = 4117.5557
$ws.Range("K73").Value = 4151.067
$ws.Range("L73").Value = 4117.5557
$ws.Range("M73").Value = -3215.067
$ws.Range("N73").Value = -5989.5557

$ws.Range("H80").Value = 2588
$ws.Range("I80").Value = 2566.6667
$ws.Range("J80").Value = 2598.6667
$ws.Range("K80").Value = 2566.6667
$ws.Range("L80").Value = 2598.6667
$ws.Range("M80").Value = -1568.6667
$ws.Range("N80").Value = -4594.6667

$ws.Range("H83").Value = 2588
$ws.Range("I83").Value = 2566.6667
$ws.Range("J83").Value = 2598.6667
$ws.Range("K83").Value = 12833.3335
$ws.Range("L83").Value = 12993.3335
$ws.Range("M83").Value = -7841.333500000001
$ws.Range("N83").Value = -22977.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 544.1667
$ws.Range("I113").Value = 482.22223
$ws.Range("J113").Value = 730
$ws.Range("K113").Value = 1446.66669
$ws.Range("L113").Value = 2190
$ws.Range("M113").Value = 723.33331
$ws.Range("N113").Value = -6530
